$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row appended below the existing table (row 6).
# Columns: A=ENTITY ID, B=MIGRATION DATE, C=ADDRESS, D=FINANCIAL INSTITUTION NAME, E=CITY
$ws.Cells.Item(6, 1).Value = "555EEE555"

# Column B holds dates stored as plain text (e.g. "2025-10-20" in B2:B5),
# so force Text formatting first to stop Excel from auto-converting the
# string into a date serial number, then drop back to the default style
# so the cell ends up unstyled like its siblings.
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 2).Value = "2025-10-25"
$ws.Cells.Item(6, 2).Style = "Normal"

$ws.Cells.Item(6, 3).Value = "T Nagar"
$ws.Cells.Item(6, 4).Value = "EEE"
$ws.Cells.Item(6, 5).Value = ""
